$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text before writing, so numeric-looking values
# (e.g. "247.55") are stored as text like the source data, not auto-converted
# to numbers by Excel's smart entry.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "36.994.61"
$ws.Range("D3").Value = "2.041.02"
$ws.Range("D4").Value = "1.01"
$ws.Range("D5").Value = "247.55"
$ws.Range("D6").Value = "0.661"
$ws.Range("D7").Value = "58.33"
$ws.Range("D9").Value = "0.381"
$ws.Range("D10").Value = "0.0780"
$ws.Range("D12").Value = "15.67"
$ws.Range("D13").Value = "2.338.76"
$ws.Range("D14").Value = "0.824"
$ws.Range("D15").Value = "5.70"
$ws.Range("D16").Value = "2.038.43"
$ws.Range("D17").Value = "18.11"
$ws.Range("D18").Value = "37.052.98"
$ws.Range("D19").Value = "74.41"
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("D21").Value = "5.30"
$ws.Range("D22").Value = "235.62"
$ws.Range("D24").Value = "2.42"
$ws.Range("D25").Value = "168.37"
$ws.Range("D26").Value = "2.15"
$ws.Range("D27").Value = "9.29"
$ws.Range("D28").Value = "19.89"
$ws.Range("D29").Value = "0.124"
$ws.Range("D30").Value = "1.12"
$ws.Range("D31").Value = "4.72"
$ws.Range("D32").Value = "0.0620"
$ws.Range("D33").Value = "4.47"
$ws.Range("D34").Value = "0.0891"
$ws.Range("D36").Value = "2.19"
$ws.Range("D37").Value = "1.76"
$ws.Range("D38").Value = "1.33"
$ws.Range("D39").Value = "0.105"
$ws.Range("D41").Value = "5.17"
$ws.Range("D42").Value = "0.0220"
$ws.Range("D43").Value = "17.25"
$ws.Range("D44").Value = "1.12"
$ws.Range("D45").Value = "95.12"
$ws.Range("D46").Value = "2.44"
$ws.Range("D47").Value = "2.90"
$ws.Range("D48").Value = "1.272.34"
$ws.Range("D49").Value = "6.75"
$ws.Range("D50").Value = "2.218.96"
$ws.Range("D51").Value = "43.37"

# Restore the original (unstyled) formatting now that the text is committed.
$dRange.ClearFormats()

# Plain text / percentage cells need no special handling.
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  +6.86%  "
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("E17").Value = "  +25.65%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E26").Value = "  +7.28%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  -3.74%  "
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("E40").Value = "  +14.08%  "
$ws.Range("E41").Value = "  +15.97%  "
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("E43").Value = "  -5.32%  "
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -0.46%  "
